$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Cells.Item(2, 8).Value = 1644.2858
$ws.Cells.Item(2, 9).Value = 3400
$ws.Cells.Item(2, 10).Value = 327.5
$ws.Cells.Item(2, 11).Value = 3400
$ws.Cells.Item(2, 12).Value = 327.5
$ws.Cells.Item(2, 13).Value = -3287
$ws.Cells.Item(2, 14).Value = -553.5

# Row 28
$ws.Cells.Item(28, 8).Value = 211.5
$ws.Cells.Item(28, 9).Value = 213.33333
$ws.Cells.Item(28, 10).Value = 195
$ws.Cells.Item(28, 11).Value = 213.33333
$ws.Cells.Item(28, 12).Value = 195
$ws.Cells.Item(28, 13).Value = 271.66667
$ws.Cells.Item(28, 14).Value = -1165

# Row 45
$ws.Cells.Item(45, 8).Value = 20251
$ws.Cells.Item(45, 10).Value = 20251
$ws.Cells.Item(45, 12).Value = 60753
$ws.Cells.Item(45, 14).Value = -61137

# Row 70
$ws.Cells.Item(70, 8).Value = 1228.7646
$ws.Cells.Item(70, 9).Value = 974.5714
$ws.Cells.Item(70, 10).Value = 1538.2174
$ws.Cells.Item(70, 11).Value = 2923.7142
$ws.Cells.Item(70, 12).Value = 4614.6522
$ws.Cells.Item(70, 13).Value = -2653.7142
$ws.Cells.Item(70, 14).Value = -5154.6522

# Row 73
$ws.Cells.Item(73, 8).Value = 1228.7646
$ws.Cells.Item(73, 9).Value = 974.5714
$ws.Cells.Item(73, 10).Value = 1538.2174
$ws.Cells.Item(73, 11).Value = 2923.7142
$ws.Cells.Item(73, 12).Value = 4614.6522
$ws.Cells.Item(73, 13).Value = -1987.7142
$ws.Cells.Item(73, 14).Value = -6486.6522

# Row 80
$ws.Cells.Item(80, 8).Value = 373.4
$ws.Cells.Item(80, 9).Value = 378.6875
$ws.Cells.Item(80, 10).Value = 352.25
$ws.Cells.Item(80, 11).Value = 1136.0625
$ws.Cells.Item(80, 12).Value = 1056.75
$ws.Cells.Item(80, 13).Value = -138.0625
$ws.Cells.Item(80, 14).Value = -3052.75

# Row 83
$ws.Cells.Item(83, 8).Value = 373.4
$ws.Cells.Item(83, 9).Value = 378.6875
$ws.Cells.Item(83, 10).Value = 352.25
$ws.Cells.Item(83, 11).Value = 3408.1875
$ws.Cells.Item(83, 12).Value = 3170.25
$ws.Cells.Item(83, 13).Value = 1583.8125
$ws.Cells.Item(83, 14).Value = -13154.25

# Row 116
$ws.Cells.Item(116, 8).Value = 5058.8184
$ws.Cells.Item(116, 9).Value = 2235.2856
$ws.Cells.Item(116, 11).Value = 2235.2856
$ws.Cells.Item(116, 13).Value = 1206.7144

# Row 125
$ws.Cells.Item(125, 8).Value = 7048.727
$ws.Cells.Item(125, 10).Value = 7504
$ws.Cells.Item(125, 12).Value = 67536
$ws.Cells.Item(125, 14).Value = -72456

# Row 135
$ws.Cells.Item(135, 8).Value = 4208.067
$ws.Cells.Item(135, 9).Value = 252.1
$ws.Cells.Item(135, 10).Value = 12120
$ws.Cells.Item(135, 11).Value = 2268.9
$ws.Cells.Item(135, 12).Value = 109080
$ws.Cells.Item(135, 13).Value = 266.0999999999999
$ws.Cells.Item(135, 14).Value = -114150

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 5220551
$ws.Cells.Item(32, 9).Value = 7379.8936
$ws.Cells.Item(32, 10).Value = 19633436
$ws.Cells.Item(32, 11).Value = 7379.8936
$ws.Cells.Item(32, 12).Value = 19633436
$ws.Cells.Item(32, 13).Value = -7092.8936
$ws.Cells.Item(32, 14).Value = -19634010

# Row 45
$ws.Cells.Item(45, 8).Value = 3096.5
$ws.Cells.Item(45, 9).Value = 1819.6364
$ws.Cells.Item(45, 10).Value = 4032.8667
$ws.Cells.Item(45, 11).Value = 1819.6364
$ws.Cells.Item(45, 12).Value = 4032.8667
$ws.Cells.Item(45, 13).Value = -1442.6364
$ws.Cells.Item(45, 14).Value = -4786.8667

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Cells.Item(7, 8).Value = 7601.5713
$ws.Cells.Item(7, 9).Value = 11138
$ws.Cells.Item(7, 11).Value = 11138
$ws.Cells.Item(7, 13).Value = -11025

# Row 94
$ws.Cells.Item(94, 8).Value = 1568.6
$ws.Cells.Item(94, 9).Value = 1669.4286
$ws.Cells.Item(94, 10).Value = 1514.3077
$ws.Cells.Item(94, 11).Value = 1669.4286
$ws.Cells.Item(94, 12).Value = 1514.3077
$ws.Cells.Item(94, 13).Value = -1218.4286
$ws.Cells.Item(94, 14).Value = -2416.3077

# Row 99
$ws.Cells.Item(99, 8).Value = 2530.7932
$ws.Cells.Item(99, 9).Value = 2064
$ws.Cells.Item(99, 10).Value = 2816.0557
$ws.Cells.Item(99, 11).Value = 2064
$ws.Cells.Item(99, 12).Value = 2816.0557
$ws.Cells.Item(99, 13).Value = -566
$ws.Cells.Item(99, 14).Value = -5812.0557

# Row 107
$ws.Cells.Item(107, 8).Value = 1033.8334
$ws.Cells.Item(107, 9).Value = 545.6667
$ws.Cells.Item(107, 10).Value = 2498.3333
$ws.Cells.Item(107, 11).Value = 545.6667
$ws.Cells.Item(107, 12).Value = 2498.3333
$ws.Cells.Item(107, 13).Value = 1374.3333
$ws.Cells.Item(107, 14).Value = -6338.3333

# Row 126
$ws.Cells.Item(126, 8).Value = 2530.7932
$ws.Cells.Item(126, 9).Value = 2064
$ws.Cells.Item(126, 10).Value = 2816.0557
$ws.Cells.Item(126, 11).Value = 6192
$ws.Cells.Item(126, 12).Value = 8448.167099999999
$ws.Cells.Item(126, 13).Value = -3722
$ws.Cells.Item(126, 14).Value = -13388.1671

$ws = $wb.Worksheets.Item("CUL")
# Row 75
$ws.Cells.Item(75, 8).Value = 1589.375
$ws.Cells.Item(75, 9).Value = 900
$ws.Cells.Item(75, 10).Value = 1819.1666
$ws.Cells.Item(75, 11).Value = 2700
$ws.Cells.Item(75, 12).Value = 5457.4998
$ws.Cells.Item(75, 13).Value = -1702
$ws.Cells.Item(75, 14).Value = -7453.4998

# Row 78
$ws.Cells.Item(78, 8).Value = 1589.375
$ws.Cells.Item(78, 9).Value = 900
$ws.Cells.Item(78, 10).Value = 1819.1666
$ws.Cells.Item(78, 11).Value = 8100
$ws.Cells.Item(78, 12).Value = 16372.4994
$ws.Cells.Item(78, 13).Value = -3108
$ws.Cells.Item(78, 14).Value = -26356.4994

# Row 98
$ws.Cells.Item(98, 8).Value = 18717.334
$ws.Cells.Item(98, 9).Value = 533.3333
$ws.Cells.Item(98, 10).Value = 36901.332
$ws.Cells.Item(98, 11).Value = 1599.9999
$ws.Cells.Item(98, 12).Value = 110703.996
$ws.Cells.Item(98, 13).Value = -101.9999
$ws.Cells.Item(98, 14).Value = -113699.996

# Row 109
$ws.Cells.Item(109, 8).Value = 4684.6665
$ws.Cells.Item(109, 9).Value = 4684.6665
$ws.Cells.Item(109, 10).Value = 0
$ws.Cells.Item(109, 11).Value = 14053.9995
$ws.Cells.Item(109, 12).Value = 0
$ws.Cells.Item(109, 13).Value = -13013.9995
$ws.Cells.Item(109, 14).ClearContents()

# Row 113
$ws.Cells.Item(113, 8).Value = 545
$ws.Cells.Item(113, 9).Value = 501.5
$ws.Cells.Item(113, 10).Value = 603
$ws.Cells.Item(113, 11).Value = 1504.5
$ws.Cells.Item(113, 12).Value = 1809
$ws.Cells.Item(113, 13).Value = 665.5
$ws.Cells.Item(113, 14).Value = -6149

# Row 129
$ws.Cells.Item(129, 8).Value = 2062.0588
$ws.Cells.Item(129, 9).Value = 390
$ws.Cells.Item(129, 10).Value = 3943.125
$ws.Cells.Item(129, 11).Value = 1170
$ws.Cells.Item(129, 12).Value = 11829.375
$ws.Cells.Item(129, 13).Value = 3830
$ws.Cells.Item(129, 14).Value = -21829.375

# Row 131
$ws.Cells.Item(131, 8).Value = 1283120.4
$ws.Cells.Item(131, 9).Value = 453.8889
$ws.Cells.Item(131, 10).Value = 5900720
$ws.Cells.Item(131, 11).Value = 1361.6667
$ws.Cells.Item(131, 12).Value = 17702160
$ws.Cells.Item(131, 13).Value = 3678.3333
$ws.Cells.Item(131, 14).Value = -17712240

$ws = $wb.Worksheets.Item("LTW")
# Row 36
$ws.Cells.Item(36, 8).Value = 34995.5
$ws.Cells.Item(36, 10).Value = 34995.5
$ws.Cells.Item(36, 12).Value = 34995.5
$ws.Cells.Item(36, 14).Value = -36119.5

# Row 40
$ws.Cells.Item(40, 8).Value = 4306.8125
$ws.Cells.Item(40, 9).Value = 2173.0908
$ws.Cells.Item(40, 11).Value = 2173.0908
$ws.Cells.Item(40, 13).Value = -2037.0908

$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Cells.Item(2, 8).Value = 260035000
$ws.Cells.Item(2, 10).Value = 333380000
$ws.Cells.Item(2, 12).Value = 333380000
$ws.Cells.Item(2, 14).Value = -333380224
